$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 194.25
$ws.Range("I5").Value = 194.25
$ws.Range("K5").Value = 194.25
$ws.Range("M5").Value = -79.25
$ws.Range("H9").Value = 140
$ws.Range("I9").Value = 146.2
$ws.Range("J9").Value = 124.5
$ws.Range("K9").Value = 146.2
$ws.Range("L9").Value = 124.5
$ws.Range("M9").Value = 22.80000000000001
$ws.Range("N9").Value = -462.5
$ws.Range("H43").Value = 4599.5
$ws.Range("I43").Value = 4900
$ws.Range("J43").Value = 4299
$ws.Range("K43").Value = 4900
$ws.Range("L43").Value = 4299
$ws.Range("M43").Value = -4831
$ws.Range("N43").Value = -4437
$ws.Range("H138").Value = 3032.1292
$ws.Range("J138").Value = 3024.6875
$ws.Range("L138").Value = 9074.0625
$ws.Range("N138").Value = -19354.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10520.625
$ws.Range("I61").Value = 9951.333000000001
$ws.Range("J61").Value = 11252.571
$ws.Range("K61").Value = 9951.333000000001
$ws.Range("L61").Value = 11252.571
$ws.Range("M61").Value = -9739.333000000001
$ws.Range("N61").Value = -11676.571
$ws.Range("H102").Value = 5242.6
$ws.Range("I102").Value = 5303.375
$ws.Range("K102").Value = 5303.375
$ws.Range("M102").Value = -3681.375
$ws.Range("H132").Value = 5700
$ws.Range("I132").Value = 3169.697
$ws.Range("K132").Value = 9509.091
$ws.Range("M132").Value = -6979.091
$ws.Range("H136").Value = 10520.625
$ws.Range("I136").Value = 9951.333000000001
$ws.Range("J136").Value = 11252.571
$ws.Range("K136").Value = 29853.999
$ws.Range("L136").Value = 33757.713
$ws.Range("M136").Value = -27303.999
$ws.Range("N136").Value = -38857.713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9666.333000000001
$ws.Range("I75").Value = 9666.333000000001
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 9666.333000000001
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -8730.333000000001
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 9666.333000000001
$ws.Range("I78").Value = 9666.333000000001
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 28998.999
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -24318.999
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 1003638.9
$ws.Range("I86").Value = 1669650.4
$ws.Range("J86").Value = 4621.75
$ws.Range("K86").Value = 1669650.4
$ws.Range("L86").Value = 4621.75
$ws.Range("M86").Value = -1668527.4
$ws.Range("N86").Value = -6867.75
$ws.Range("H89").Value = 1003638.9
$ws.Range("I89").Value = 1669650.4
$ws.Range("J89").Value = 4621.75
$ws.Range("K89").Value = 8348252
$ws.Range("L89").Value = 23108.75
$ws.Range("M89").Value = -8342636
$ws.Range("N89").Value = -34340.75
$ws.Range("H94").Value = 467.375
$ws.Range("J94").Value = 699.5
$ws.Range("L94").Value = 699.5
$ws.Range("N94").Value = -1601.5
$ws.Range("H99").Value = 2901
$ws.Range("I99").Value = 2901
$ws.Range("K99").Value = 2901
$ws.Range("M99").Value = -1403
$ws.Range("H105").Value = 3453.158
$ws.Range("I105").Value = 3800
$ws.Range("J105").Value = 2482
$ws.Range("K105").Value = 3800
$ws.Range("L105").Value = 2482
$ws.Range("M105").Value = -2053
$ws.Range("N105").Value = -5976

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6818.727
$ws.Range("I16").Value = 6666.3335
$ws.Range("K16").Value = 6666.3335
$ws.Range("M16").Value = -6379.3335
$ws.Range("H58").Value = 4193
$ws.Range("I58").Value = 2125.3845
$ws.Range("J58").Value = 9568.799999999999
$ws.Range("K58").Value = 2125.3845
$ws.Range("L58").Value = 9568.799999999999
$ws.Range("M58").Value = -1922.3845
$ws.Range("N58").Value = -9974.799999999999
$ws.Range("H107").Value = 878.0625
$ws.Range("I107").Value = 681.8946999999999
$ws.Range("K107").Value = 681.8946999999999
$ws.Range("M107").Value = 1238.1053
$ws.Range("H113").Value = 6818.727
$ws.Range("I113").Value = 6666.3335
$ws.Range("K113").Value = 6666.3335
$ws.Range("M113").Value = -4496.3335
$ws.Range("H136").Value = 4193
$ws.Range("I136").Value = 2125.3845
$ws.Range("J136").Value = 9568.799999999999
$ws.Range("K136").Value = 6376.1535
$ws.Range("L136").Value = 28706.4
$ws.Range("M136").Value = -3826.1535
$ws.Range("N136").Value = -33806.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 739.6
$ws.Range("I117").Value = 649.5
$ws.Range("J117").Value = 799.6667
$ws.Range("K117").Value = 1948.5
$ws.Range("L117").Value = 2399.0001
$ws.Range("M117").Value = 1493.5
$ws.Range("N117").Value = -9283.000100000001
$ws.Range("H137").Value = 2802.3076
$ws.Range("I137").Value = 2048.111
$ws.Range("K137").Value = 6144.333
$ws.Range("M137").Value = -1044.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16705
$ws.Range("J15").Value = 16705
$ws.Range("L15").Value = 16705
$ws.Range("N15").Value = -17281
$ws.Range("H70").Value = 9197.299999999999
$ws.Range("J70").Value = 10747.25
$ws.Range("L70").Value = 10747.25
$ws.Range("N70").Value = -11287.25
$ws.Range("H73").Value = 9197.299999999999
$ws.Range("J73").Value = 10747.25
$ws.Range("L73").Value = 10747.25
$ws.Range("N73").Value = -12619.25
$ws.Range("H81").Value = 16705
$ws.Range("J81").Value = 16705
$ws.Range("L81").Value = 16705
$ws.Range("N81").Value = -18701
$ws.Range("H84").Value = 16705
$ws.Range("J84").Value = 16705
$ws.Range("L84").Value = 50115
$ws.Range("N84").Value = -60099
$ws.Range("H92").Value = 31070.143
$ws.Range("J92").Value = 31070.143
$ws.Range("L92").Value = 31070.143
$ws.Range("N92").Value = -34814.143
$ws.Range("H107").Value = 428.42856
$ws.Range("I107").Value = 294.875
$ws.Range("K107").Value = 294.875
$ws.Range("M107").Value = 1625.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2768.6
$ws.Range("I22").Value = 2318
$ws.Range("J22").Value = 3219.2
$ws.Range("K22").Value = 2318
$ws.Range("L22").Value = 3219.2
$ws.Range("M22").Value = -2023
$ws.Range("N22").Value = -3809.2
$ws.Range("H27").Value = 2768.6
$ws.Range("I27").Value = 2318
$ws.Range("J27").Value = 3219.2
$ws.Range("K27").Value = 2318
$ws.Range("L27").Value = 3219.2
$ws.Range("M27").Value = -2211
$ws.Range("N27").Value = -3433.2
$ws.Range("H100").Value = 5016.591
$ws.Range("I100").Value = 5345.2354
$ws.Range("K100").Value = 5345.2354
$ws.Range("M100").Value = -4804.2354
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2599.5454
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2599.5454
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -40608
$ws.Range("H122").Value = 4816.5386
$ws.Range("I122").Value = 3172.8096
$ws.Range("J122").Value = 11720.2
$ws.Range("K122").Value = 9518.4288
$ws.Range("L122").Value = 35160.60000000001
$ws.Range("M122").Value = -7068.4288
$ws.Range("N122").Value = -40060.60000000001
$ws.Range("H132").Value = 18404.654
$ws.Range("I132").Value = 11299.405
$ws.Range("J132").Value = 30923.428
$ws.Range("K132").Value = 33898.215
$ws.Range("L132").Value = 92770.284
$ws.Range("M132").Value = -31368.215
$ws.Range("N132").Value = -97830.284
